$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 'ECs'
$ws.Cells.Item(2, 2).Value = 'Il34'
$ws.Cells.Item(2, 3).Value = 'Csf1r'
$ws.Cells.Item(2, 4).Value = 'ECs'
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 2.11382
$ws.Cells.Item(2, 8).Value = 6.34146
$ws.Cells.Item(2, 9).Value = 0.1611092823235492
$ws.Cells.Item(2, 10).Value = 0.1611092823235492
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.9601406666666668
$ws.Cells.Item(2, 14).Value = 2.880422
$ws.Cells.Item(2, 15).Value = 0.004534252661098308
$ws.Cells.Item(2, 16).Value = 0.004534252661098308
$ws.Cells.Item(2, 17).Value = 2.029564544013334
$ws.Cells.Item(2, 18).Value = 18.26608089612
$ws.Cells.Item(2, 19).Value = 0.0007305101921031915
$ws.Cells.Item(2, 20).Value = 0.0007305101921031913

# Row 3
$ws.Cells.Item(3, 1).Value = 'ECs'
$ws.Cells.Item(3, 2).Value = 'Il34'
$ws.Cells.Item(3, 3).Value = 'Csf1r'
$ws.Cells.Item(3, 4).Value = 'FAPs'
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 2.11382
$ws.Cells.Item(3, 8).Value = 6.34146
$ws.Cells.Item(3, 9).Value = 0.1611092823235492
$ws.Cells.Item(3, 10).Value = 0.1611092823235492
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 3.928236333333333
$ws.Cells.Item(3, 14).Value = 11.784709
$ws.Cells.Item(3, 15).Value = 0.01855104847259158
$ws.Cells.Item(3, 16).Value = 0.01855104847259158
$ws.Cells.Item(3, 17).Value = 8.303584526126667
$ws.Cells.Item(3, 18).Value = 74.73226073513999
$ws.Cells.Item(3, 19).Value = 0.002988746105768602
$ws.Cells.Item(3, 20).Value = 0.002988746105768602

# Row 4
$ws.Cells.Item(4, 1).Value = 'ECs'
$ws.Cells.Item(4, 2).Value = 'Il34'
$ws.Cells.Item(4, 3).Value = 'Csf1r'
$ws.Cells.Item(4, 4).Value = 'MuSCs'
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 2.11382
$ws.Cells.Item(4, 8).Value = 6.34146
$ws.Cells.Item(4, 9).Value = 0.1611092823235492
$ws.Cells.Item(4, 10).Value = 0.1611092823235492
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1.316544333333334
$ws.Cells.Item(4, 14).Value = 3.949633
$ws.Cells.Item(4, 15).Value = 0.006217364657196653
$ws.Cells.Item(4, 16).Value = 0.006217364657196652
$ws.Cells.Item(4, 17).Value = 2.782937742686667
$ws.Cells.Item(4, 18).Value = 25.04643968418
$ws.Cells.Item(4, 19).Value = 0.001001675157864752
$ws.Cells.Item(4, 20).Value = 0.001001675157864752

# Row 5
$ws.Cells.Item(5, 1).Value = 'ECs'
$ws.Cells.Item(5, 2).Value = 'Il34'
$ws.Cells.Item(5, 3).Value = 'Csf1r'
$ws.Cells.Item(5, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 2.11382
$ws.Cells.Item(5, 8).Value = 6.34146
$ws.Cells.Item(5, 9).Value = 0.1611092823235492
$ws.Cells.Item(5, 10).Value = 0.1611092823235492
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 205.547872
$ws.Cells.Item(5, 14).Value = 616.643616
$ws.Cells.Item(5, 15).Value = 0.9706973342091134
$ws.Cells.Item(5, 16).Value = 0.9706973342091134
$ws.Cells.Item(5, 17).Value = 434.49120279104
$ws.Cells.Item(5, 18).Value = 3910.42082511936
$ws.Cells.Item(5, 19).Value = 0.1563883508678126
$ws.Cells.Item(5, 20).Value = 0.1563883508678126

# Row 6
$ws.Cells.Item(6, 1).Value = 'FAPs'
$ws.Cells.Item(6, 2).Value = 'Il34'
$ws.Cells.Item(6, 3).Value = 'Csf1r'
$ws.Cells.Item(6, 4).Value = 'ECs'
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 4.352037
$ws.Cells.Item(6, 8).Value = 13.056111
$ws.Cells.Item(6, 9).Value = 0.3316997462960574
$ws.Cells.Item(6, 10).Value = 0.3316997462960574
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.9601406666666668
$ws.Cells.Item(6, 14).Value = 2.880422
$ws.Cells.Item(6, 15).Value = 0.004534252661098308
$ws.Cells.Item(6, 16).Value = 0.004534252661098308
$ws.Cells.Item(6, 17).Value = 4.178567706538001
$ws.Cells.Item(6, 18).Value = 37.607109358842
$ws.Cells.Item(6, 19).Value = 0.001504010457328532
$ws.Cells.Item(6, 20).Value = 0.001504010457328532

# Row 7
$ws.Cells.Item(7, 1).Value = 'FAPs'
$ws.Cells.Item(7, 2).Value = 'Il34'
$ws.Cells.Item(7, 3).Value = 'Csf1r'
$ws.Cells.Item(7, 4).Value = 'FAPs'
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 4.352037
$ws.Cells.Item(7, 8).Value = 13.056111
$ws.Cells.Item(7, 9).Value = 0.3316997462960574
$ws.Cells.Item(7, 10).Value = 0.3316997462960574
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 3.928236333333333
$ws.Cells.Item(7, 14).Value = 11.784709
$ws.Cells.Item(7, 15).Value = 0.01855104847259158
$ws.Cells.Item(7, 16).Value = 0.01855104847259158
$ws.Cells.Item(7, 17).Value = 17.095829867411
$ws.Cells.Item(7, 18).Value = 153.862468806699
$ws.Cells.Item(7, 19).Value = 0.006153378071884489
$ws.Cells.Item(7, 20).Value = 0.006153378071884489

# Row 8
$ws.Cells.Item(8, 1).Value = 'FAPs'
$ws.Cells.Item(8, 2).Value = 'Il34'
$ws.Cells.Item(8, 3).Value = 'Csf1r'
$ws.Cells.Item(8, 4).Value = 'MuSCs'
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 4.352037
$ws.Cells.Item(8, 8).Value = 13.056111
$ws.Cells.Item(8, 9).Value = 0.3316997462960574
$ws.Cells.Item(8, 10).Value = 0.3316997462960574
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 1.316544333333334
$ws.Cells.Item(8, 14).Value = 3.949633
$ws.Cells.Item(8, 15).Value = 0.006217364657196653
$ws.Cells.Item(8, 16).Value = 0.006217364657196652
$ws.Cells.Item(8, 17).Value = 5.729649650807001
$ws.Cells.Item(8, 18).Value = 51.566846857263
$ws.Cells.Item(8, 19).Value = 0.002062298279422203
$ws.Cells.Item(8, 20).Value = 0.002062298279422203

# Row 9
$ws.Cells.Item(9, 1).Value = 'FAPs'
$ws.Cells.Item(9, 2).Value = 'Il34'
$ws.Cells.Item(9, 3).Value = 'Csf1r'
$ws.Cells.Item(9, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 4.352037
$ws.Cells.Item(9, 8).Value = 13.056111
$ws.Cells.Item(9, 9).Value = 0.3316997462960574
$ws.Cells.Item(9, 10).Value = 0.3316997462960574
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 205.547872
$ws.Cells.Item(9, 14).Value = 616.643616
$ws.Cells.Item(9, 15).Value = 0.9706973342091134
$ws.Cells.Item(9, 16).Value = 0.9706973342091134
$ws.Cells.Item(9, 17).Value = 894.5519442152639
$ws.Cells.Item(9, 18).Value = 8050.967497937375
$ws.Cells.Item(9, 19).Value = 0.3219800594874221
$ws.Cells.Item(9, 20).Value = 0.3219800594874221

# Row 10
$ws.Cells.Item(10, 1).Value = 'MuSCs'
$ws.Cells.Item(10, 2).Value = 'Il34'
$ws.Cells.Item(10, 3).Value = 'Csf1r'
$ws.Cells.Item(10, 4).Value = 'ECs'
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 6.620393
$ws.Cells.Item(10, 8).Value = 19.861179
$ws.Cells.Item(10, 9).Value = 0.5045873181869075
$ws.Cells.Item(10, 10).Value = 0.5045873181869074
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.9601406666666668
$ws.Cells.Item(10, 14).Value = 2.880422
$ws.Cells.Item(10, 15).Value = 0.004534252661098308
$ws.Cells.Item(10, 16).Value = 0.004534252661098308
$ws.Cells.Item(10, 17).Value = 6.356508548615334
$ws.Cells.Item(10, 18).Value = 57.20857693753801
$ws.Cells.Item(10, 19).Value = 0.002287926390245444
$ws.Cells.Item(10, 20).Value = 0.002287926390245444

# Row 11
$ws.Cells.Item(11, 1).Value = 'MuSCs'
$ws.Cells.Item(11, 2).Value = 'Il34'
$ws.Cells.Item(11, 3).Value = 'Csf1r'
$ws.Cells.Item(11, 4).Value = 'FAPs'
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 6.620393
$ws.Cells.Item(11, 8).Value = 19.861179
$ws.Cells.Item(11, 9).Value = 0.5045873181869075
$ws.Cells.Item(11, 10).Value = 0.5045873181869074
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 3.928236333333333
$ws.Cells.Item(11, 14).Value = 11.784709
$ws.Cells.Item(11, 15).Value = 0.01855104847259158
$ws.Cells.Item(11, 16).Value = 0.01855104847259158
$ws.Cells.Item(11, 17).Value = 26.00646832354566
$ws.Cells.Item(11, 18).Value = 234.058214911911
$ws.Cells.Item(11, 19).Value = 0.00936062379834031
$ws.Cells.Item(11, 20).Value = 0.009360623798340309

# Row 12
$ws.Cells.Item(12, 1).Value = 'MuSCs'
$ws.Cells.Item(12, 2).Value = 'Il34'
$ws.Cells.Item(12, 3).Value = 'Csf1r'
$ws.Cells.Item(12, 4).Value = 'MuSCs'
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 6.620393
$ws.Cells.Item(12, 8).Value = 19.861179
$ws.Cells.Item(12, 9).Value = 0.5045873181869075
$ws.Cells.Item(12, 10).Value = 0.5045873181869074
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 1.316544333333334
$ws.Cells.Item(12, 14).Value = 3.949633
$ws.Cells.Item(12, 15).Value = 0.006217364657196653
$ws.Cells.Item(12, 16).Value = 0.006217364657196652
$ws.Cells.Item(12, 17).Value = 8.716040888589667
$ws.Cells.Item(12, 18).Value = 78.44436799730701
$ws.Cells.Item(12, 19).Value = 0.00313720335856492
$ws.Cells.Item(12, 20).Value = 0.003137203358564919

# Row 13
$ws.Cells.Item(13, 1).Value = 'MuSCs'
$ws.Cells.Item(13, 2).Value = 'Il34'
$ws.Cells.Item(13, 3).Value = 'Csf1r'
$ws.Cells.Item(13, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 6.620393
$ws.Cells.Item(13, 8).Value = 19.861179
$ws.Cells.Item(13, 9).Value = 0.5045873181869075
$ws.Cells.Item(13, 10).Value = 0.5045873181869074
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 205.547872
$ws.Cells.Item(13, 14).Value = 616.643616
$ws.Cells.Item(13, 15).Value = 0.9706973342091134
$ws.Cells.Item(13, 16).Value = 0.9706973342091134
$ws.Cells.Item(13, 17).Value = 1360.807692953696
$ws.Cells.Item(13, 18).Value = 12247.26923658326
$ws.Cells.Item(13, 19).Value = 0.4898015646397568
$ws.Cells.Item(13, 20).Value = 0.4898015646397567

# Row 14
$ws.Cells.Item(14, 1).Value = 'Resolving-Mac'
$ws.Cells.Item(14, 2).Value = 'Il34'
$ws.Cells.Item(14, 3).Value = 'Csf1r'
$ws.Cells.Item(14, 4).Value = 'ECs'
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.034161
$ws.Cells.Item(14, 8).Value = 0.102483
$ws.Cells.Item(14, 9).Value = 0.002603653193486089
$ws.Cells.Item(14, 10).Value = 0.002603653193486088
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.9601406666666668
$ws.Cells.Item(14, 14).Value = 2.880422
$ws.Cells.Item(14, 15).Value = 0.004534252661098308
$ws.Cells.Item(14, 16).Value = 0.004534252661098308
$ws.Cells.Item(14, 17).Value = 0.032799365314
$ws.Cells.Item(14, 18).Value = 0.295194287826
$ws.Cells.Item(14, 19).Value = 0.00001180562142114141
$ws.Cells.Item(14, 20).Value = 0.0000118056214211414

# Row 15
$ws.Cells.Item(15, 1).Value = 'Resolving-Mac'
$ws.Cells.Item(15, 2).Value = 'Il34'
$ws.Cells.Item(15, 3).Value = 'Csf1r'
$ws.Cells.Item(15, 4).Value = 'FAPs'
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.034161
$ws.Cells.Item(15, 8).Value = 0.102483
$ws.Cells.Item(15, 9).Value = 0.002603653193486089
$ws.Cells.Item(15, 10).Value = 0.002603653193486088
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 3.928236333333333
$ws.Cells.Item(15, 14).Value = 11.784709
$ws.Cells.Item(15, 15).Value = 0.01855104847259158
$ws.Cells.Item(15, 16).Value = 0.01855104847259158
$ws.Cells.Item(15, 17).Value = 0.134192481383
$ws.Cells.Item(15, 18).Value = 1.207732332447
$ws.Cells.Item(15, 19).Value = 0.00004830049659817829
$ws.Cells.Item(15, 20).Value = 0.00004830049659817828

# Row 16
$ws.Cells.Item(16, 1).Value = 'Resolving-Mac'
$ws.Cells.Item(16, 2).Value = 'Il34'
$ws.Cells.Item(16, 3).Value = 'Csf1r'
$ws.Cells.Item(16, 4).Value = 'MuSCs'
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.034161
$ws.Cells.Item(16, 8).Value = 0.102483
$ws.Cells.Item(16, 9).Value = 0.002603653193486089
$ws.Cells.Item(16, 10).Value = 0.002603653193486088
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 1.316544333333334
$ws.Cells.Item(16, 14).Value = 3.949633
$ws.Cells.Item(16, 15).Value = 0.006217364657196653
$ws.Cells.Item(16, 16).Value = 0.006217364657196652
$ws.Cells.Item(16, 17).Value = 0.04497447097100001
$ws.Cells.Item(16, 18).Value = 0.404770238739
$ws.Cells.Item(16, 19).Value = 0.00001618786134477761
$ws.Cells.Item(16, 20).Value = 0.0000161878613447776

# Row 17
$ws.Cells.Item(17, 1).Value = 'Resolving-Mac'
$ws.Cells.Item(17, 2).Value = 'Il34'
$ws.Cells.Item(17, 3).Value = 'Csf1r'
$ws.Cells.Item(17, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(17, 5).Value = 1
$ws.Cells.Item(17, 6).Value = 0.3333333333333333
$ws.Cells.Item(17, 7).Value = 0.034161
$ws.Cells.Item(17, 8).Value = 0.102483
$ws.Cells.Item(17, 9).Value = 0.002603653193486089
$ws.Cells.Item(17, 10).Value = 0.002603653193486088
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 205.547872
$ws.Cells.Item(17, 14).Value = 616.643616
$ws.Cells.Item(17, 15).Value = 0.9706973342091134
$ws.Cells.Item(17, 16).Value = 0.9706973342091134
$ws.Cells.Item(17, 17).Value = 7.021720855392
$ws.Cells.Item(17, 18).Value = 63.195487698528
$ws.Cells.Item(17, 19).Value = 0.002527359214121991
$ws.Cells.Item(17, 20).Value = 0.002527359214121991
